$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (STRASSE, HAUSNR) before the existing PLZ column (J),
# pushing PLZ/WOHNORT/VERMÖGEN/HASEL/HASSH two columns to the right.
$ws.Range("J1:K1").EntireColumn.Insert()

# New column widths (matching the bestFit width used by the neighbouring
# AUFENTHALTSBEWILLIGUNG column, which is stored as 23.1640625 points; the
# closest value reachable through the ColumnWidth (character-unit) setter
# is used here since it is internally rounded to whole pixels).
$ws.Columns("J:K").ColumnWidth = 22.3

# Street + house-number data, filled row by row (matching the order the
# original author entered the values in Excel), then the header labels.
$ws.Range("J2").Value = "Ackerstrasse"
$ws.Range("K2").Value = 11

$ws.Range("K1").Value = "HAUSNR"
$ws.Range("J1").Value = "STRASSE"

$ws.Range("J3").Value = "Denzingsteig"
$ws.Range("K3").Value = 22

$ws.Range("J4").Value = "Eichenweg"
$ws.Range("K4").Value = 33

# Selection ends up on the newly added cells.
$ws.Range("J3:K3").Select()
